$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H48").Value = 10342.286
$ws.Range("J48").Value = 10342.286
$ws.Range("L48").Value = 31026.858
$ws.Range("N48").Value = -31610.858
$ws.Range("H56").Value = 10342.286
$ws.Range("J56").Value = 10342.286
$ws.Range("L56").Value = 31026.858
$ws.Range("N56").Value = -32094.858
$ws.Range("H127").Value = 3302.077
$ws.Range("I127").Value = 1532.8572
$ws.Range("J127").Value = 5366.1665
$ws.Range("K127").Value = 4598.571599999999
$ws.Range("L127").Value = 16098.4995
$ws.Range("M127").Value = 361.4284000000007
$ws.Range("N127").Value = -26018.4995
$ws.Range("H129").Value = 901.2
$ws.Range("I129").Value = 901.2
$ws.Range("J129").Value = 0
$ws.Range("K129").Value = 2703.6
$ws.Range("L129").Value = 0
$ws.Range("M129").Value = 2296.4
$ws.Range("N129").ClearContents()
$ws.Range("H131").Value = 3409.8572
$ws.Range("I131").Value = 3311.5
$ws.Range("K131").Value = 9934.5
$ws.Range("M131").Value = -4894.5
$ws.Range("H138").Value = 3589.3062
$ws.Range("I138").Value = 3089.8845
$ws.Range("J138").Value = 4153.8696
$ws.Range("K138").Value = 9269.6535
$ws.Range("L138").Value = 12461.6088
$ws.Range("M138").Value = -4129.6535
$ws.Range("N138").Value = -22741.6088
$ws.Range("H141").Value = 2063.0908
$ws.Range("I141").Value = 1747.6538
$ws.Range("J141").Value = 3234.7144
$ws.Range("K141").Value = 5242.9614
$ws.Range("L141").Value = 9704.143199999999
$ws.Range("M141").Value = -62.96140000000014
$ws.Range("N141").Value = -20064.1432

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6173.2324
$ws.Range("I32").Value = 4551.422
$ws.Range("J32").Value = 14586.375
$ws.Range("K32").Value = 4551.422
$ws.Range("L32").Value = 14586.375
$ws.Range("M32").Value = -4264.422
$ws.Range("N32").Value = -15160.375
$ws.Range("H132").Value = 3958.487
$ws.Range("I132").Value = 1754.4849
$ws.Range("K132").Value = 5263.4547
$ws.Range("M132").Value = -2733.4547
$ws.Range("H139").Value = 118966.664
$ws.Range("J139").Value = 118966.664
$ws.Range("L139").Value = 118966.664
$ws.Range("N139").Value = -129246.664

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7217.385
$ws.Range("I31").Value = 2671.3438
$ws.Range("K31").Value = 2671.3438
$ws.Range("M31").Value = -2376.3438
$ws.Range("H34").Value = 7217.385
$ws.Range("I34").Value = 2671.3438
$ws.Range("K34").Value = 2671.3438
$ws.Range("M34").Value = -2469.3438
$ws.Range("H122").Value = 2708.9333
$ws.Range("I122").Value = 2247.5
$ws.Range("J122").Value = 3977.875
$ws.Range("K122").Value = 6742.5
$ws.Range("L122").Value = 11933.625
$ws.Range("M122").Value = -4292.5
$ws.Range("N122").Value = -16833.625

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 9482932
$ws.Range("I4").Value = 20810972
$ws.Range("J4").Value = 42899.125
$ws.Range("K4").Value = 62432916
$ws.Range("L4").Value = 128697.375
$ws.Range("M4").Value = -62432804
$ws.Range("N4").Value = -128921.375
$ws.Range("H129").Value = 11908237
$ws.Range("J129").Value = 20837600
$ws.Range("L129").Value = 62512800
$ws.Range("N129").Value = -62522800
$ws.Range("H139").Value = 4475.8066
$ws.Range("I139").Value = 3182.2
$ws.Range("K139").Value = 9546.599999999999
$ws.Range("M139").Value = -4406.599999999999
$ws.Range("H140").Value = 1701.9117
$ws.Range("I140").Value = 1528.0769
$ws.Range("J140").Value = 1809.5238
$ws.Range("K140").Value = 4584.2307
$ws.Range("L140").Value = 5428.5714
$ws.Range("M140").Value = 595.7692999999999
$ws.Range("N140").Value = -15788.5714

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H134").Value = 311331.25
$ws.Range("J134").Value = 311331.25
$ws.Range("L134").Value = 933993.75
$ws.Range("N134").Value = -939063.75
$ws.Range("H136").Value = 74999
$ws.Range("J136").Value = 74999
$ws.Range("L136").Value = 224997
$ws.Range("N136").Value = -230097

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3209.6843
$ws.Range("I22").Value = 3483.5386
$ws.Range("K22").Value = 3483.5386
$ws.Range("M22").Value = -3188.5386
$ws.Range("H27").Value = 3209.6843
$ws.Range("I27").Value = 3483.5386
$ws.Range("K27").Value = 3483.5386
$ws.Range("M27").Value = -3376.5386
$ws.Range("H40").Value = 5656.4443
$ws.Range("I40").Value = 4484.6665
$ws.Range("K40").Value = 4484.6665
$ws.Range("M40").Value = -4348.6665
$ws.Range("H46").Value = 2915.1667
$ws.Range("I46").Value = 2064.3333
$ws.Range("J46").Value = 3766
$ws.Range("K46").Value = 2064.3333
$ws.Range("L46").Value = 3766
$ws.Range("M46").Value = -1876.3333
$ws.Range("N46").Value = -4142
$ws.Range("H68").Value = 9880.652
$ws.Range("I68").Value = 7733.1816
$ws.Range("J68").Value = 11849.167
$ws.Range("K68").Value = 7733.1816
$ws.Range("L68").Value = 11849.167
$ws.Range("M68").Value = -6984.1816
$ws.Range("N68").Value = -13347.167
$ws.Range("H71").Value = 9880.652
$ws.Range("I71").Value = 7733.1816
$ws.Range("J71").Value = 11849.167
$ws.Range("K71").Value = 38665.908
$ws.Range("L71").Value = 59245.835
$ws.Range("M71").Value = -34921.908
$ws.Range("N71").Value = -66733.83499999999
$ws.Range("H131").Value = 79081.5
$ws.Range("J131").Value = 79081.5
$ws.Range("L131").Value = 79081.5
$ws.Range("N131").Value = -89161.5
$ws.Range("H136").Value = 10551.107
$ws.Range("I136").Value = 9711.25
$ws.Range("J136").Value = 10691.083
$ws.Range("K136").Value = 29133.75
$ws.Range("L136").Value = 32073.249
$ws.Range("M136").Value = -26583.75
$ws.Range("N136").Value = -37173.249

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 8399
$ws.Range("I62").Value = 3995
$ws.Range("J62").Value = 9500
$ws.Range("K62").Value = 3995
$ws.Range("L62").Value = 9500
$ws.Range("M62").Value = -3371
$ws.Range("N62").Value = -10748
$ws.Range("H65").Value = 8399
$ws.Range("I65").Value = 3995
$ws.Range("J65").Value = 9500
$ws.Range("K65").Value = 19975
$ws.Range("L65").Value = 47500
$ws.Range("M65").Value = -16855
$ws.Range("N65").Value = -53740
$ws.Range("H122").Value = 2920.68
$ws.Range("I122").Value = 2924.0212
$ws.Range("J122").Value = 2868.3333
$ws.Range("K122").Value = 8772.063600000001
$ws.Range("L122").Value = 8604.999899999999
$ws.Range("M122").Value = -6322.063600000001
$ws.Range("N122").Value = -13504.9999
$ws.Range("H136").Value = 9965.1
$ws.Range("I136").Value = 14750.333
$ws.Range("J136").Value = 7914.2856
$ws.Range("K136").Value = 44250.999
$ws.Range("L136").Value = 23742.8568
$ws.Range("M136").Value = -41700.999
$ws.Range("N136").Value = -28842.8568
